# Neues csv Template -> Schild Rosette Drücker hinzugefügt
# Appends new leader-log rows (209-216) to the log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2025-10-14T16:59:53", "25.0803_Akon_Dermanence_Grundrissneu_ab_dwg.3dm", "bbf7e392-3a37-4155-b244-b1382392e1d0", "schiebetuere", "Standard 1:10 Schiebetürbeschriftung", "schiebetuere.csv"),
    @("2025-10-16T09:33:05", "Untitled", "5f187c27-5d2a-4f8f-9ac3-288c2cbbe980", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuere.csv"),
    @("2025-10-16T09:33:34", "Untitled", "8a4bf1bf-617c-42a6-8cc3-a45bdc7d3eb6", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuere.csv"),
    @("2025-10-16T09:36:20", "Untitled", "986f81b8-ab95-4f8b-8a90-e8004634d8d2", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuere.csv"),
    @("2025-10-16T09:36:36", "Untitled", "7c73a897-bb99-457f-97d1-b3483bc29f2a", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuere.csv"),
    @("2025-10-16T09:36:57", "Untitled", "0ce0b9b6-3457-4d1b-aabc-a99c490daff3", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuere.csv"),
    @("2025-10-16T09:37:20", "Untitled", "2bdc1dc8-b6ee-4642-b318-ffd17e70545e", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuere.csv"),
    @("2025-10-16T09:37:35", "Untitled", "f6f1397b-917d-4dc7-b4ff-6e5329e2330c", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuere.csv")
)

$startRow = 209
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]

    # normalize to the sheet's default (unformatted) style, same as the
    # other data rows above, instead of inheriting the blank-row style
    $ws.Range("A" + $r + ":F" + $r).VerticalAlignment = -4107
}
